# Commit: "built the word docs"
#
# The rebuilt document no longer emits the (unused) built-in
# "Footnote Text" paragraph style definition
# (<w:style w:type="paragraph" w:styleId="FootnoteText"> ... </w:style>)
# in word/styles.xml - basedOn Normal, w:next pointing back at itself,
# uiPriority 9, unhideWhenUsed/qFormat flags, no overrides. It isn't
# referenced by any paragraph or footnote in the document, so removing
# the style definition is the whole edit.
$d = $word.ActiveDocument

try {
    $style = $d.Styles("FootnoteText")
    if ($style -ne $null) {
        $style.Delete()
    }
} catch {
    # Style already absent - nothing to do.
}
